$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bahrain / Bahreïn -> add trailing asterisk (row 73: BHR)
$ws.Range("B73").Value = "Bahrain*"
$ws.Range("C73").Value = "Bahreïn*"

# Cabo Verde -> add trailing asterisk (row 95: CPV)
$ws.Range("B95").Value = "Cabo Verde*"
$ws.Range("C95").Value = "Cabo Verde*"

# Nigeria -> remove trailing asterisk (row 182: NGA)
$ws.Range("B182").Value = "Nigeria"
$ws.Range("C182").Value = "Nigeria"

# Papua New Guinea / Papouasie-Nouvelle-Guinée -> add trailing asterisk (row 195: PNG)
$ws.Range("B195").Value = "Papua New Guinea*"
$ws.Range("C195").Value = "Papouasie-Nouvelle-Guinée*"

# South Sudan / Soudan du Sud -> add trailing asterisk (row 214: SSD)
$ws.Range("B214").Value = "South Sudan*"
$ws.Range("C214").Value = "Soudan du Sud*"
